# Applies the "chore: update Sheets via scheduled runner" commit:
# recalculated currentAveragePrice* / LevePrice* / LeveProfit* figures
# (columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve tables.
$wb = $excel.ActiveWorkbook

# hunk 0 sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 120
$ws.Range("I4").Value = 120
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 120
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -6

# hunk 1 sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4409.25
$ws.Range("I74").Value = 3933.3333
$ws.Range("J74").Value = 4694.8
$ws.Range("K74").Value = 3933.3333
$ws.Range("L74").Value = 4694.8
$ws.Range("M74").Value = -2997.3333
$ws.Range("N74").Value = -6566.8

# hunk 2 sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4409.25
$ws.Range("I77").Value = 3933.3333
$ws.Range("J77").Value = 4694.8
$ws.Range("K77").Value = 19666.6665
$ws.Range("L77").Value = 23474
$ws.Range("M77").Value = -14986.6665
$ws.Range("N77").Value = -32834

# hunk 3 sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 88237120
$ws.Range("I86").Value = 214287090
$ws.Range("J86").Value = 2139.8
$ws.Range("K86").Value = 214287090
$ws.Range("L86").Value = 2139.8
$ws.Range("M86").Value = -214285967
$ws.Range("N86").Value = -4385.8

# hunk 4 sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 88237120
$ws.Range("I89").Value = 214287090
$ws.Range("J89").Value = 2139.8
$ws.Range("K89").Value = 1071435450
$ws.Range("L89").Value = 10699
$ws.Range("M89").Value = -1071429834
$ws.Range("N89").Value = -21931

# hunk 5 sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2490.6316
$ws.Range("I137").Value = 2107.2354
$ws.Range("J137").Value = 5749.5
$ws.Range("K137").Value = 6321.706200000001
$ws.Range("L137").Value = 17248.5
$ws.Range("M137").Value = -3771.706200000001
$ws.Range("N137").Value = -22348.5

# hunk 6 sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1977.1466
$ws.Range("I138").Value = 1509
$ws.Range("J138").Value = 2289.2444
$ws.Range("K138").Value = 4527
$ws.Range("L138").Value = 6867.733200000001
$ws.Range("M138").Value = 613
$ws.Range("N138").Value = -17147.7332

# hunk 7 sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()

# hunk 8 sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 898.9
$ws.Range("I74").Value = 817.6667
$ws.Range("J74").Value = 1142.6
$ws.Range("K74").Value = 817.6667
$ws.Range("L74").Value = 1142.6
$ws.Range("M74").Value = 56.33330000000001
$ws.Range("N74").Value = -2890.6

# hunk 9 sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 898.9
$ws.Range("I77").Value = 817.6667
$ws.Range("J77").Value = 1142.6
$ws.Range("K77").Value = 4088.3335
$ws.Range("L77").Value = 5713
$ws.Range("M77").Value = 279.6665000000003
$ws.Range("N77").Value = -14449

# hunk 10 sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 34939.5
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 34939.5
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 34939.5
$ws.Range("N128").Value = -44899.5

# hunk 11 sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H123").Value = 31992.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 31992.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 31992.5
$ws.Range("N123").Value = -41792.5

# hunk 12 sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2500750
$ws.Range("I6").Value = 5000000
$ws.Range("J6").Value = 1500
$ws.Range("K6").Value = 5000000
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -4999887
$ws.Range("N6").Value = -1726

# hunk 13 sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 144.66667
$ws.Range("I7").Value = 110
$ws.Range("J7").Value = 188
$ws.Range("K7").Value = 110
$ws.Range("L7").Value = 188
$ws.Range("M7").Value = 3
$ws.Range("N7").Value = -414

# hunk 14 sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()

# hunk 15 sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 27175
$ws.Range("I12").Value = 2700
$ws.Range("J12").Value = 35333.332
$ws.Range("K12").Value = 2700
$ws.Range("L12").Value = 35333.332
$ws.Range("M12").Value = -2530
$ws.Range("N12").Value = -35673.332

# hunk 16 sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 32500
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 32500
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 32500
$ws.Range("N120").Value = -39758

# hunk 17 sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 36284.168
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 36284.168
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 36284.168
$ws.Range("N133").Value = -41344.168

# hunk 18 sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 189285.14
$ws.Range("I141").Value = 200000
$ws.Range("J141").Value = 187499.33
$ws.Range("K141").Value = 200000
$ws.Range("L141").Value = 187499.33
$ws.Range("M141").Value = -194820
$ws.Range("N141").Value = -197859.33

# hunk 19 sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 9098.583000000001
$ws.Range("I122").Value = 714.8570999999999
$ws.Range("J122").Value = 20835.8
$ws.Range("K122").Value = 6433.7139
$ws.Range("L122").Value = 187522.2
$ws.Range("M122").Value = -3983.7139
$ws.Range("N122").Value = -192422.2

# hunk 20 sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2249.2122
$ws.Range("I132").Value = 1722.7273
$ws.Range("J132").Value = 2512.4546
$ws.Range("K132").Value = 15504.5457
$ws.Range("L132").Value = 22612.0914
$ws.Range("M132").Value = -12974.5457
$ws.Range("N132").Value = -27672.0914

# hunk 21 sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 13902914
$ws.Range("I137").Value = 55593892
$ws.Range("J137").Value = 5922.5557
$ws.Range("K137").Value = 166781676
$ws.Range("L137").Value = 17767.6671
$ws.Range("M137").Value = -166776576
$ws.Range("N137").Value = -27967.6671

# hunk 22 sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2075.3845
$ws.Range("I140").Value = 1580
$ws.Range("J140").Value = 2500
$ws.Range("K140").Value = 4740
$ws.Range("L140").Value = 7500
$ws.Range("M140").Value = 440
$ws.Range("N140").Value = -17860

# hunk 23 sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 56668
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 56668
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 56668
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -56948

# hunk 24 sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 7000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 7000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 7000
$ws.Range("N18").Value = -7586

# hunk 25 sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 6900
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 6900
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 6900
$ws.Range("N21").Value = -7246

# hunk 26 sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 6900
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 6900
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 6900
$ws.Range("N30").Value = -7110

# hunk 27 sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 54499.25
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 54499.25
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 54499.25
$ws.Range("N112").Value = -56715.25

# hunk 28 sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H120").Value = 34050
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 34050
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 34050
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -43726

# hunk 29 sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 29800
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 29800
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 29800
$ws.Range("N131").Value = -39880

# hunk 30 sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 70007
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 70007
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 70007
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -70467

# hunk 31 sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 939814.5600000001
$ws.Range("I25").Value = 7585.6
$ws.Range("J25").Value = 1716672
$ws.Range("K25").Value = 7585.6
$ws.Range("L25").Value = 1716672
$ws.Range("M25").Value = -7355.6
$ws.Range("N25").Value = -1717132

# hunk 32 sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 28500
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 28500
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 28500
$ws.Range("N109").Value = -31274

# hunk 33 sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H131").Value = 30000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 30000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080

# hunk 34 sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2169.0322
$ws.Range("I132").Value = 1514.4286
$ws.Range("J132").Value = 3543.7
$ws.Range("K132").Value = 4543.2858
$ws.Range("L132").Value = 10631.1
$ws.Range("M132").Value = -2013.2858
$ws.Range("N132").Value = -15691.1

# hunk 35 sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 3255250
$ws.Range("I139").Value = 19300000
$ws.Range("J139").Value = 46300
$ws.Range("K139").Value = 19300000
$ws.Range("L139").Value = 46300
$ws.Range("M139").Value = -19294860
$ws.Range("N139").Value = -56580

# hunk 36 sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 23108.572
$ws.Range("I123").Value = 20000
$ws.Range("J123").Value = 24352
$ws.Range("K123").Value = 20000
$ws.Range("L123").Value = 24352
$ws.Range("M123").Value = -15100
$ws.Range("N123").Value = -34152

